$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append new row 28 with the new mail-log entry ---
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A28").Value = "Klacht over levering"
$ws.Range("B28").Value = "mailmind.test@zohomail.eu"
$ws.Range("C28").Value = "Ik ben niet tevreden over mijn bestelling. Ik hoor graag hoe jullie dit oplossen."
$ws.Range("D28").Value = "Klacht / Probleem"
$ws.Range("F28").Value = "2025-06-19 21:45:10"
$ws.Range("G28").Value = "Nee"

# Extend the conditional formatting ranges to cover the newly added row 28
$dCond = $ws.Range("D2:D27").FormatConditions
$dCond.Item(1).ModifyAppliesToRange($ws.Range("D2:D28"))

$gCond = $ws.Range("G2:G27").FormatConditions
$gCond.Item(1).ModifyAppliesToRange($ws.Range("G2:G28"))

# --- Sheet "Dashboard": re-order the category counts now that
#     "Klacht / Probleem" moved from 1 to 2 occurrences ---
$ws2 = $wb.Worksheets.Item("Dashboard")

$ws2.Range("A7").Value = "Klacht / Probleem"
$ws2.Range("B7").Value = 2

$ws2.Range("A8").Value = "Offerte / Prijsaanvraag"
$ws2.Range("B8").Value = 2

$ws2.Range("A9").Value = "Openingstijden / Locatie"
$ws2.Range("B9").Value = 1

$ws2.Range("A10").Value = "Sollicitatie / Vacature"
$ws2.Range("B10").Value = 1
